$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.900.61'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.639.25'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').Value = '''217.12'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '''0.505'
$ws.Range('E6').Value = '  +2.04%  '
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').Value = '''0.253'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').Value = '''0.0624'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  +3.36%  '
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.868.09'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = '1.656.28'
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('E14').Value = '  -1.17%  '
$ws.Range('D15').Value = '''0.531'
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('D16').Value = '''67.08'
$ws.Range('E16').Value = '  +2.90%  '
$ws.Range('D17').Value = '26.892.41'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '''218.20'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('D21').Value = '''6.72'
$ws.Range('E21').Value = '  +2.34%  '
$ws.Range('D22').Value = '''4.40'
$ws.Range('E22').Value = '  +0.78%  '
$ws.Range('E23').Value = '  +2.46%  '
$ws.Range('D24').Value = '''9.17'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '''147.40'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('D29').Value = '''15.73'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').Value = '''3.33'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').Value = '''3.00'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D35').Value = '1.264.50'
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('D38').Value = '''0.838'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').Value = '''0.810'
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('D43').Value = '1.779.31'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').Value = '''62.18'
$ws.Range('E44').Value = '  +1.84%  '
$ws.Range('D45').Value = '''2.12'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('D46').Value = '''92.17'
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +9.09%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.0512'
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.67'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.0960'
$ws.Range('E51').Value = '  -0.82%  '
